$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Наличные")

$data = @(
    ,@(7676096317, 'Ходатайство студенты 2025.docx', 3, 0.6, '2025-07-11 21:49:13')
    ,@(7676096317, 'Ходатайство студенты 2025.docx', 3, 0.6, '2025-07-11 21:54:47')
    ,@(7676096317, 'Ходатайство студенты 2025.docx', 3, 0.6, '2025-07-11 21:56:23')
    ,@(7676096317, 'Ходатайство студенты 2025.docx', 3, 0.6, '2025-07-11 22:00:02')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 22:01:28')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 22:15:49')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 22:23:18')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 22:25:03')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 22:34:32')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 22:38:22')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 22:58:58')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:00:11')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:03:27')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:11:38')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:17:10')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:34:21')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:36:44')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:41:03')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:44:26')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:45:42')
    ,@(7676096317, 'отчет маз 9docx.pdf', 40, 8, '2025-07-11 23:49:01')
)

$startRow = 68
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
